$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in newly-discovered worker IDs for three existing rows ---
# Joshua Ladd (row 19)
$ws.Range("C19").Value = "CB46D992"
# Kelly Pum (row 21)
$ws.Range("C21").Value = "11D14592"
# Michael Makris (row 24)
$ws.Range("C24").Value = "D4C6CF96"

# --- Add a new worker row (38): An Nguyen ---
$ws.Range("A38").Value = "An Nguyen"
$ws.Range("B38").Value = 1218
$ws.Range("C38").Value = "AF8446CB"

# Row 37 (ZyAsia Holmes) carries a highlighted border style. Touch the
# border formatting on the two rows below it (38-39) so row 39 stays part
# of the sheet's used range (matching the trailing blank formatted row in
# the source) without inheriting that highlight.
$ws.Range("A38:A39").Borders.LineStyle = -4142

# --- View state: scroll/selection left where the editor was last working ---
$ws.Range("G21").Select()
